$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1124
$ws.Range("I62").Value = 1144.6666
$ws.Range("K62").Value = 1144.6666
$ws.Range("M62").Value = -520.6666
$ws.Range("H65").Value = 1124
$ws.Range("I65").Value = 1144.6666
$ws.Range("K65").Value = 5723.333000000001
$ws.Range("M65").Value = -2603.333000000001
$ws.Range("H86").Value = 1111.7273
$ws.Range("I86").Value = 1080.7778
$ws.Range("K86").Value = 1080.7778
$ws.Range("M86").Value = 42.22219999999993
$ws.Range("H89").Value = 1111.7273
$ws.Range("I89").Value = 1080.7778
$ws.Range("K89").Value = 5403.889
$ws.Range("M89").Value = 212.1109999999999
$ws.Range("H112").Value = 2256.0605
$ws.Range("J112").Value = 2340.3547
$ws.Range("L112").Value = 7021.0641
$ws.Range("N112").Value = -9237.0641
$ws.Range("H129").Value = 1092.7667
$ws.Range("J129").Value = 1248
$ws.Range("L129").Value = 3744
$ws.Range("N129").Value = -13744
$ws.Range("H132").Value = 860.1087
$ws.Range("I132").Value = 803.4878
$ws.Range("K132").Value = 2410.4634
$ws.Range("M132").Value = 119.5365999999999
$ws.Range("H137").Value = 1918
$ws.Range("H141").Value = 7002913
$ws.Range("I141").Value = 14000548
$ws.Range("K141").Value = 42001644
$ws.Range("M141").Value = -41996464

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6078.742
$ws.Range("I32").Value = 4023.7817
$ws.Range("J32").Value = 16353.546
$ws.Range("K32").Value = 4023.7817
$ws.Range("L32").Value = 16353.546
$ws.Range("M32").Value = -3736.7817
$ws.Range("N32").Value = -16927.546
$ws.Range("H61").Value = 5658.4
$ws.Range("I61").Value = 5334.364
$ws.Range("J61").Value = 6549.5
$ws.Range("K61").Value = 5334.364
$ws.Range("L61").Value = 6549.5
$ws.Range("M61").Value = -5122.364
$ws.Range("N61").Value = -6973.5
$ws.Range("H74").Value = 1003.4828
$ws.Range("I74").Value = 604.2917
$ws.Range("K74").Value = 604.2917
$ws.Range("M74").Value = 269.7083
$ws.Range("H77").Value = 1003.4828
$ws.Range("I77").Value = 604.2917
$ws.Range("K77").Value = 3021.4585
$ws.Range("M77").Value = 1346.5415
$ws.Range("H136").Value = 5658.4
$ws.Range("I136").Value = 5334.364
$ws.Range("J136").Value = 6549.5
$ws.Range("K136").Value = 16003.092
$ws.Range("L136").Value = 19648.5
$ws.Range("M136").Value = -13453.092
$ws.Range("N136").Value = -24748.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 158086.16
$ws.Range("I86").Value = 4593.3335
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 4593.3335
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = -3470.3335
$ws.Range("N86").Value = -2002246
$ws.Range("H89").Value = 158086.16
$ws.Range("I89").Value = 4593.3335
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 22966.6675
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = -17350.6675
$ws.Range("N89").Value = -10011232
$ws.Range("H105").Value = 2378.5
$ws.Range("I105").Value = 2167
$ws.Range("K105").Value = 2167
$ws.Range("M105").Value = -420

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3890.5715
$ws.Range("I31").Value = 5095.3335
$ws.Range("K31").Value = 5095.3335
$ws.Range("M31").Value = -4800.3335
$ws.Range("H34").Value = 3890.5715
$ws.Range("I34").Value = 5095.3335
$ws.Range("K34").Value = 5095.3335
$ws.Range("M34").Value = -4893.3335
$ws.Range("H99").Value = 2514
$ws.Range("I99").Value = 2473.3333
$ws.Range("J99").Value = 2562.8
$ws.Range("K99").Value = 2473.3333
$ws.Range("L99").Value = 2562.8
$ws.Range("M99").Value = -975.3332999999998
$ws.Range("N99").Value = -5558.8
$ws.Range("H126").Value = 2514
$ws.Range("I126").Value = 2473.3333
$ws.Range("J126").Value = 2562.8
$ws.Range("K126").Value = 7419.999899999999
$ws.Range("L126").Value = 7688.400000000001
$ws.Range("M126").Value = -4949.999899999999
$ws.Range("N126").Value = -12628.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 100011800
$ws.Range("I17").Value = 500000200
$ws.Range("J17").Value = 14700.5
$ws.Range("K17").Value = 1500000600
$ws.Range("L17").Value = 44101.5
$ws.Range("M17").Value = -1500000431
$ws.Range("N17").Value = -44439.5
$ws.Range("H60").Value = 269.5
$ws.Range("I60").Value = 269
$ws.Range("K60").Value = 807
$ws.Range("M60").Value = -556
$ws.Range("H126").Value = 2876.6667
$ws.Range("I126").Value = 1315
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 3945
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = 995
$ws.Range("N126").Value = -27880
$ws.Range("H131").Value = 15005.9795
$ws.Range("J131").Value = 15951.392
$ws.Range("L131").Value = 47854.176
$ws.Range("N131").Value = -57934.176
$ws.Range("H133").Value = 4388.8887
$ws.Range("I133").Value = 3333.3333
$ws.Range("K133").Value = 9999.999899999999
$ws.Range("M133").Value = -4939.999899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 4569.6665
$ws.Range("I70").Value = 4133
$ws.Range("K70").Value = 4133
$ws.Range("M70").Value = -3863
$ws.Range("H73").Value = 4569.6665
$ws.Range("I73").Value = 4133
$ws.Range("K73").Value = 4133
$ws.Range("M73").Value = -3197
$ws.Range("H80").Value = 1980.8182
$ws.Range("I80").Value = 1754.1428
$ws.Range("J80").Value = 2377.5
$ws.Range("K80").Value = 1754.1428
$ws.Range("L80").Value = 2377.5
$ws.Range("M80").Value = -756.1428000000001
$ws.Range("N80").Value = -4373.5
$ws.Range("H83").Value = 1980.8182
$ws.Range("I83").Value = 1754.1428
$ws.Range("J83").Value = 2377.5
$ws.Range("K83").Value = 8770.714
$ws.Range("L83").Value = 11887.5
$ws.Range("M83").Value = -3778.714
$ws.Range("N83").Value = -21871.5
$ws.Range("H97").Value = 865.1177
$ws.Range("I97").Value = 932.4231
$ws.Range("K97").Value = 932.4231
$ws.Range("M97").Value = -436.4231
$ws.Range("H126").Value = 3773793
$ws.Range("I126").Value = 6176578
$ws.Range("K126").Value = 18529734
$ws.Range("M126").Value = -18527264
$ws.Range("H132").Value = 2567428.8
$ws.Range("I132").Value = 3498966.5
$ws.Range("J132").Value = 5699.5
$ws.Range("K132").Value = 10496899.5
$ws.Range("L132").Value = 17098.5
$ws.Range("M132").Value = -10494369.5
$ws.Range("N132").Value = -22158.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15323.917
$ws.Range("J40").Value = 12618.6
$ws.Range("L40").Value = 12618.6
$ws.Range("N40").Value = -12890.6
$ws.Range("H136").Value = 2598.8262
$ws.Range("I136").Value = 1782.9166
$ws.Range("J136").Value = 3488.9092
$ws.Range("K136").Value = 5348.7498
$ws.Range("L136").Value = 10466.7276
$ws.Range("M136").Value = -2798.7498
$ws.Range("N136").Value = -15566.7276

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 58280.355
$ws.Range("I122").Value = 100126.5
$ws.Range("K122").Value = 300379.5
$ws.Range("M122").Value = -297929.5
$ws.Range("H136").Value = 16341233
$ws.Range("I136").Value = 25253552
$ws.Range("J136").Value = 1980.25
$ws.Range("K136").Value = 75760656
$ws.Range("L136").Value = 5940.75
$ws.Range("M136").Value = -75758106
$ws.Range("N136").Value = -11040.75
